$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.538331317889367
$ws.Range("L2").Value = 0.62055734756295

$ws.Range("B3").Value = 0.458387576622686
$ws.Range("L3").Value = 0.703121504911223

$ws.Range("B4").Value = 0.320249242216214
$ws.Range("L4").Value = 0.699570192636595

$ws.Range("B5").Value = 0.319190975705175
$ws.Range("L5").Value = 0.377306157800078

$ws.Range("B6").Value = 0.010244201702362
$ws.Range("C6").Value = 0.0839217421719148
$ws.Range("L6").Value = -0.0145192493656787
